$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E values per row (row 2..6); G is recomputed as the sum of B+C+D+E
$data = @{
    2 = @(0.01514828764759746, 0.3127903958511391, 0.1575252929769615, 8.660232485948974)
    3 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974)
    4 = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732)
    5 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 8.660232485948974)
    6 = @(0.6753301551942219, 1.667794583268128, 26.21740644021617, 645.3272768299601)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = $b + $c + $d + $e
}
